# Insert a new weekly record at row 115 for "Vega Central Mapocho de Santiago"
# (Arveja Verde). This pushes the existing rows 115-143 down to 116-144,
# growing the used range from A1:R143 to A1:R144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(115).Insert()

$ws.Cells.Item(115, 1).Value = 9
$ws.Cells.Item(115, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(115, 3).Value = "Metropolitana"
$ws.Cells.Item(115, 4).Value = 44855
$ws.Cells.Item(115, 5).Value = 13
$ws.Cells.Item(115, 6).Value = 100112022
$ws.Cells.Item(115, 7).Value = "Arveja Verde"
$ws.Cells.Item(115, 8).Value = "Perfection"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 55
$ws.Cells.Item(115, 11).Value = 18000
$ws.Cells.Item(115, 12).Value = 20000
$ws.Cells.Item(115, 13).Value = 18909
$ws.Cells.Item(115, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(115, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(115, 16).Value = 756
$ws.Cells.Item(115, 17).Value = 25
$ws.Cells.Item(115, 18).Value = "Hortaliza"
